$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H = userCount. Fix a few data-import errors from jornada 6 (J6).
$ws.Cells.Item(2, 8).Value = 50
$ws.Cells.Item(3, 8).Value = 324
$ws.Cells.Item(4, 8).Value = 133
$ws.Cells.Item(6, 8).Value = 83
$ws.Cells.Item(7, 8).Value = 98
$ws.Cells.Item(8, 8).Value = 93
$ws.Cells.Item(12, 8).Value = 271
$ws.Cells.Item(13, 8).Value = 82
$ws.Cells.Item(14, 8).Value = 70
$ws.Cells.Item(15, 8).Value = 473
$ws.Cells.Item(16, 8).Value = 105
